$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25:C25").Copy()
$ws.Range("A29:C33").PasteSpecial(-4122)
